$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above row 29; this shifts the existing rows
# 29..140 down to 30..141, preserving all their data intact (including
# the D-column date style), matching the target diff.
$ws.Rows(29).Insert()

# Populate the newly inserted row 29 with the new record's data.
$ws.Range("A29").Value = 8
$ws.Range("B29").Value = "Terminal La Palmera de La Serena"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44677
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112040
$ws.Range("G29").Value = "Cilantro"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 3000
$ws.Range("K29").Value = 2500
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = 2750
$ws.Range("N29").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O29").Value = "Provincia del Elquí"
$ws.Range("P29").Value = 1833
$ws.Range("Q29").Value = 1.5
$ws.Range("R29").Value = "Hortaliza"
